$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.597878666666666
$ws.Range("H2").Value = 4.793635999999999
$ws.Range("I2").Value = 0.8992131381376172
$ws.Range("J2").Value = 0.8992131381376171
$ws.Range("M2").Value = 1.376848666666667
$ws.Range("N2").Value = 4.130546
$ws.Range("O2").Value = 0.1003061556015877
$ws.Range("P2").Value = 0.1003061556015877
$ws.Range("Q2").Value = 2.200037111695111
$ws.Range("R2").Value = 19.800334005256
$ws.Range("S2").Value = 0.09019661295302384
$ws.Range("T2").Value = 0.09019661295302382

# Row 3
$ws.Range("G3").Value = 1.597878666666666
$ws.Range("H3").Value = 4.793635999999999
$ws.Range("I3").Value = 0.8992131381376172
$ws.Range("J3").Value = 0.8992131381376171
$ws.Range("M3").Value = 2.848096333333333
$ws.Range("N3").Value = 8.544288999999999
$ws.Range("O3").Value = 0.207489465542554
$ws.Range("P3").Value = 0.207489465542554
$ws.Range("Q3").Value = 4.550912371644888
$ws.Range("R3").Value = 40.95821134480399
$ws.Range("S3").Value = 0.186577253441017
$ws.Range("T3").Value = 0.186577253441017

# Row 4
$ws.Range("G4").Value = 1.597878666666666
$ws.Range("H4").Value = 4.793635999999999
$ws.Range("I4").Value = 0.8992131381376172
$ws.Range("J4").Value = 0.8992131381376171
$ws.Range("M4").Value = 7.562766000000001
$ws.Range("N4").Value = 22.688298
$ws.Range("O4").Value = 0.5509624997574636
$ws.Range("P4").Value = 0.5509624997574636
$ws.Range("Q4").Value = 12.084382452392
$ws.Range("R4").Value = 108.759442071528
$ws.Range("S4").Value = 0.495432718403055
$ws.Range("T4").Value = 0.4954327184030549

# Row 5
$ws.Range("G5").Value = 1.597878666666666
$ws.Range("H5").Value = 4.793635999999999
$ws.Range("I5").Value = 0.8992131381376172
$ws.Range("J5").Value = 0.8992131381376171
$ws.Range("M5").Value = 1.938751333333333
$ws.Range("N5").Value = 5.816254000000001
$ws.Range("O5").Value = 0.1412418790983945
$ws.Range("P5").Value = 0.1412418790983945
$ws.Range("Q5").Value = 3.097889395504889
$ws.Range("R5").Value = 27.881004559544
$ws.Range("S5").Value = 0.1270065533405213
$ws.Range("T5").Value = 0.1270065533405212

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.1790956666666667
$ws.Range("H6").Value = 0.537287
$ws.Range("I6").Value = 0.1007868618623829
$ws.Range("J6").Value = 0.1007868618623829
$ws.Range("M6").Value = 1.376848666666667
$ws.Range("N6").Value = 4.130546
$ws.Range("O6").Value = 0.1003061556015877
$ws.Range("P6").Value = 0.1003061556015877
$ws.Range("Q6").Value = 0.2465876298557778
$ws.Range("R6").Value = 2.219288668702
$ws.Range("S6").Value = 0.01010954264856391
$ws.Range("T6").Value = 0.01010954264856391

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.1790956666666667
$ws.Range("H7").Value = 0.537287
$ws.Range("I7").Value = 0.1007868618623829
$ws.Range("J7").Value = 0.1007868618623829
$ws.Range("M7").Value = 2.848096333333333
$ws.Range("N7").Value = 8.544288999999999
$ws.Range("O7").Value = 0.207489465542554
$ws.Range("P7").Value = 0.207489465542554
$ws.Range("Q7").Value = 0.5100817115492221
$ws.Range("R7").Value = 4.590735403943
$ws.Range("S7").Value = 0.02091221210153706
$ws.Range("T7").Value = 0.02091221210153706

# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.1790956666666667
$ws.Range("H8").Value = 0.537287
$ws.Range("I8").Value = 0.1007868618623829
$ws.Range("J8").Value = 0.1007868618623829
$ws.Range("M8").Value = 7.562766000000001
$ws.Range("N8").Value = 22.688298
$ws.Range("O8").Value = 0.5509624997574636
$ws.Range("P8").Value = 0.5509624997574636
$ws.Range("Q8").Value = 1.354458618614
$ws.Range("R8").Value = 12.190127567526
$ws.Range("S8").Value = 0.05552978135440868
$ws.Range("T8").Value = 0.05552978135440868

# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.1790956666666667
$ws.Range("H9").Value = 0.537287
$ws.Range("I9").Value = 0.1007868618623829
$ws.Range("J9").Value = 0.1007868618623829
$ws.Range("M9").Value = 1.938751333333333
$ws.Range("N9").Value = 5.816254000000001
$ws.Range("O9").Value = 0.1412418790983945
$ws.Range("P9").Value = 0.1412418790983945
$ws.Range("Q9").Value = 0.3472219625442222
$ws.Range("R9").Value = 3.124997662898
$ws.Range("S9").Value = 0.01423532575787328
$ws.Range("T9").Value = 0.01423532575787328
